$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.644.62"
$ws.Range("E2").Value = "  +4.35%  "

$ws.Range("D3").Value = "1.604.68"
$ws.Range("E3").Value = "  +3.53%  "

$ws.Range("E4").Value = "  -0.44%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.518"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.89%  "

$ws.Range("E7").Value = "  -0.56%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "26.73"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +11.36%  "

$ws.Range("E9").Value = "  +3.47%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0600"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.97%  "

$ws.Range("E11").Value = "  +3.08%  "

$ws.Range("D12").Value = "1.833.05"
$ws.Range("E12").Value = "  +3.44%  "

$ws.Range("D13").Value = "1.612.48"
$ws.Range("E13").Value = "  +3.87%  "

$ws.Range("D14").Value = "29.677.16"
$ws.Range("E14").Value = "  +4.66%  "

$ws.Range("E15").Value = "  +3.87%  "

$ws.Range("E16").Value = "  +3.77%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "246.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.76%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.80%  "

$ws.Range("D20").Value = "0.0₃0696"
$ws.Range("E20").Value = "  +3.29%  "

$ws.Range("E21").Value = "  -0.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.99%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.05%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.24%  "

$ws.Range("E27").Value = "  +6.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.73%  "

$ws.Range("E29").Value = "  -0.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0474"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.42%  "

$ws.Range("E31").Value = "  +0.64%  "

$ws.Range("E32").Value = "  +2.74%  "

$ws.Range("D33").Value = "1.441.06"
$ws.Range("E33").Value = "  +4.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.88%  "

$ws.Range("E35").Value = "  -1.65%  "

$ws.Range("E36").Value = "  +11.06%  "

$ws.Range("E37").Value = "  +3.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.50%  "

$ws.Range("E39").Value = "  +2.85%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.536"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "55.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +27.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.97"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.801"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.96%  "

$ws.Range("E44").Value = "  -0.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0468"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.29%  "

$ws.Range("D48").Value = "1.744.32"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.70%  "

$ws.Range("E50").Value = "  -3.98%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0103"
$ws.Range("E51").Value = "  -0.05%  "
